# Applies the latest cryptos snapshot values (prices / 1h volume change,
# plus a couple of re-ranked coins) onto Sheet1, matching the upstream
# GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "61.875.62"
$ws.Range("E2").Value = "'" + "  +0.71%  "
$ws.Range("D3").Value = "'" + "3.421.25"
$ws.Range("E3").Value = "'" + "  +0.67%  "
$ws.Range("D4").Value = "'" + "0.999"
$ws.Range("E4").Value = "'" + "  -0.09%  "
$ws.Range("D5").Value = "'" + "408.85"
$ws.Range("E5").Value = "'" + "  +1.24%  "
$ws.Range("D6").Value = "'" + "128.50"
$ws.Range("E6").Value = "'" + "  -1.85%  "
$ws.Range("D7").Value = "'" + "0.632"
$ws.Range("E7").Value = "'" + "  +7.23%  "
$ws.Range("E8").Value = "'" + "  -0.07%  "
$ws.Range("E9").Value = "'" + "  +7.98%  "
$ws.Range("E10").Value = "'" + "  +8.89%  "
$ws.Range("D11").Value = "'" + "42.56"
$ws.Range("B12").Value = "'" + "Polkadot"
$ws.Range("C12").Value = "'" + "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'" + "9.08"
$ws.Range("E12").Value = "'" + "  +9.17%  "
$ws.Range("B13").Value = "'" + "TRON"
$ws.Range("C13").Value = "'" + "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'" + "0.141"
$ws.Range("E13").Value = "'" + "  +0.12%  "
$ws.Range("D14").Value = "'" + "3.963.75"
$ws.Range("E14").Value = "'" + "  +0.87%  "
$ws.Range("D15").Value = "'" + "21.22"
$ws.Range("E15").Value = "'" + "  +7.70%  "
$ws.Range("D16").Value = "'" + "0.0000203"
$ws.Range("E16").Value = "'" + "  +41.68%  "
$ws.Range("D17").Value = "'" + "3.452.88"
$ws.Range("E17").Value = "'" + "  +1.77%  "
$ws.Range("D18").Value = "'" + "12.30"
$ws.Range("E18").Value = "'" + "  +5.20%  "
$ws.Range("D19").Value = "'" + "1.08"
$ws.Range("E19").Value = "'" + "  +6.81%  "
$ws.Range("D20").Value = "'" + "61.855.23"
$ws.Range("E20").Value = "'" + "  +0.78%  "
$ws.Range("D21").Value = "'" + "445.21"
$ws.Range("E21").Value = "'" + "  +43.28%  "
$ws.Range("D22").Value = "'" + "90.99"
$ws.Range("E22").Value = "'" + "  +9.48%  "
$ws.Range("D23").Value = "'" + "3.20"
$ws.Range("E23").Value = "'" + "  +1.68%  "
$ws.Range("D24").Value = "'" + "12.93"
$ws.Range("E24").Value = "'" + "  +1.90%  "
$ws.Range("E25").Value = "'" + "  +2.97%  "
$ws.Range("D26").Value = "'" + "32.91"
$ws.Range("E26").Value = "'" + "  +11.70%  "
$ws.Range("D27").Value = "'" + "8.77"
$ws.Range("E27").Value = "'" + "  +9.41%  "
$ws.Range("E28").Value = "'" + "  -0.33%  "
$ws.Range("D29").Value = "'" + "2.76"
$ws.Range("E29").Value = "'" + "  +2.06%  "
$ws.Range("D30").Value = "'" + "7.60"
$ws.Range("E30").Value = "'" + "  -5.77%  "
$ws.Range("D31").Value = "'" + "11.92"
$ws.Range("E31").Value = "'" + "  +5.89%  "
$ws.Range("E32").Value = "'" + "  -0.08%  "
$ws.Range("E33").Value = "'" + "  -0.17%  "
$ws.Range("D34").Value = "'" + "42.62"
$ws.Range("E34").Value = "'" + "  -1.70%  "
$ws.Range("E35").Value = "'" + "  -0.02%  "
$ws.Range("D36").Value = "'" + "0.0498"
$ws.Range("E36").Value = "'" + "  +3.75%  "
$ws.Range("D37").Value = "'" + "53.21"
$ws.Range("E37").Value = "'" + "  +3.75%  "
$ws.Range("D38").Value = "'" + "0.998"
$ws.Range("E38").Value = "'" + "  +0.11%  "
$ws.Range("E39").Value = "'" + "  +1.23%  "
$ws.Range("E40").Value = "'" + "  +8.05%  "
$ws.Range("E41").Value = "'" + "  -0.60%  "
$ws.Range("E42").Value = "'" + "  -1.87%  "
$ws.Range("D43").Value = "'" + "141.91"
$ws.Range("E43").Value = "'" + "  +1.85%  "
$ws.Range("E44").Value = "'" + "  +8.28%  "
$ws.Range("E45").Value = "'" + "  +1.07%  "
$ws.Range("D46").Value = "'" + "2.45"
$ws.Range("E46").Value = "'" + "  +10.72%  "
$ws.Range("D47").Value = "'" + "16.55"
$ws.Range("E47").Value = "'" + "  +0.07%  "
$ws.Range("D48").Value = "'" + "22.32"
$ws.Range("E48").Value = "'" + "  +5.90%  "
$ws.Range("D49").Value = "'" + "3.768.06"
$ws.Range("E49").Value = "'" + "  +0.96%  "
$ws.Range("B50").Value = "'" + "ThetaToken"
$ws.Range("C50").Value = "'" + "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'" + "2.08"
$ws.Range("E50").Value = "'" + "  +9.07%  "
$ws.Range("B51").Value = "'" + "Maker"
$ws.Range("C51").Value = "'" + "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'" + "2.124.39"
$ws.Range("E51").Value = "'" + "  +1.77%  "
